# This script re-applies the upstream "Automatic update of files" edit to the
# Artfynd worksheet. The export tool re-ordered the observation rows (moving the
# two rows that were at sheet positions 8 and 20 - "Korallrot" and "Bollvitmossa" -
# up to the top, directly below the header row) and every other row shifted down
# to fill the gap, in its original relative order. Concretely this means:
#   new row 2  <- old row 8
#   new row 3  <- old row 20
#   new row 4  <- old row 2
#   new row 5  <- old row 3
#   new row 6  <- old row 4
#   new row 7  <- old row 5
#   new row 8  <- old row 6
#   new row 9  <- old row 7
#   new row 10 <- old row 9
#   new row 11 <- old row 10
#   new row 12 <- old row 11
#   new row 13 <- old row 12
#   new row 14 <- old row 13
#   new row 15 <- old row 14
#   new row 16 <- old row 15
#   new row 17 <- old row 16
#   new row 18 <- old row 17
#   new row 19 <- old row 18
#   new row 20 <- old row 19
# Rather than physically moving rows (which this engine does not reliably support
# as an atomic "insert cut cells" operation), we directly overwrite every cell that
# differs between the original row and its new content with the value coming from
# its source row, leaving all cells that keep the same value untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 101774542
$ws.Range("B2").Value = 96237
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 220093
$ws.Range("F2").Value = "Korallrot"
$ws.Range("G2").Value = "Corallorhiza trifida"
$ws.Range("H2").Value = "Châtel."
$ws.Range("I2").Value = "'8"
$ws.Range("J2").Value = "stjälkar/strån/skott"
$ws.Range("K2").Value = "blomning"
$ws.Range("Q2").Value = 567847.3860201587
$ws.Range("R2").Value = 6676956.727084515
$ws.Range("Z2").Value = "17:05"
$ws.Range("AB2").Value = "17:05"
$ws.Range("AW2").Value = "Annelie Hilmerby"
$ws.Range("AX2").Value = "Annelie Hilmerby"

# Row 3
$ws.Range("A3").Value = 101774450
$ws.Range("B3").Value = 93868
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 2869
$ws.Range("F3").Value = "Bollvitmossa"
$ws.Range("G3").Value = "Sphagnum wulfianum"
$ws.Range("H3").Value = "Girg."
$ws.Range("I3").ClearContents()
$ws.Range("Q3").Value = 567869.409834059
$ws.Range("R3").Value = 6676977.014183999
$ws.Range("Z3").Value = "17:01"
$ws.Range("AB3").Value = "17:01"
$ws.Range("AW3").Value = "Philipp Weiss"
$ws.Range("AX3").Value = "Philipp Weiss"

# Row 4
$ws.Range("A4").Value = 101775819
$ws.Range("B4").Value = 96334
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = "'3"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("Q4").Value = 567890.146411219
$ws.Range("R4").Value = 6676825.326947801
$ws.Range("Z4").Value = "17:55"
$ws.Range("AB4").Value = "17:55"
$ws.Range("AC4").ClearContents()

# Row 5
$ws.Range("A5").Value = 101775508
$ws.Range("I5").Value = "'10"
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("Q5").Value = 567931.7906619203
$ws.Range("R5").Value = 6676834.54900873
$ws.Range("Z5").Value = "17:37"
$ws.Range("AB5").Value = "17:37"
$ws.Range("AW5").Value = "fanny westling"
$ws.Range("AX5").Value = "fanny westling"

# Row 6
$ws.Range("A6").Value = 101774746
$ws.Range("B6").Value = 56411
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100049
$ws.Range("F6").Value = "Spillkråka"
$ws.Range("G6").Value = "Dryocopus martius"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("Q6").Value = 567863.9375606392
$ws.Range("R6").Value = 6676842.731599597
$ws.Range("Z6").Value = "17:17"
$ws.Range("AB6").Value = "17:17"
$ws.Range("AC6").Value = "Mat hål"
$ws.Range("AW6").Value = "FREDRIK  Månsson "
$ws.Range("AX6").Value = "FREDRIK  Månsson "

# Row 7
$ws.Range("A7").Value = 101775559
$ws.Range("I7").Value = "'40"
$ws.Range("K7").Value = "fullt utvecklade blad"
$ws.Range("Q7").Value = 567920.9165920488
$ws.Range("R7").Value = 6676830.371183628
$ws.Range("Z7").Value = "17:41"
$ws.Range("AB7").Value = "17:41"
$ws.Range("AW7").Value = "Annelie Hilmerby"
$ws.Range("AX7").Value = "Annelie Hilmerby"

# Row 8
$ws.Range("A8").Value = 101779914
$ws.Range("B8").Value = 96334
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("I8").Value = "'1"
$ws.Range("J8").Value = "plantor/tuvor"
$ws.Range("K8").Value = "fullt utvecklade blad"
$ws.Range("Q8").Value = 567971.7663531044
$ws.Range("R8").Value = 6676853.183049097
$ws.Range("Z8").Value = "20:43"
$ws.Range("AB8").Value = "20:43"

# Row 9
$ws.Range("A9").Value = 101775856
$ws.Range("I9").Value = "'3"
$ws.Range("Q9").Value = 567863.1820311426
$ws.Range("R9").Value = 6676829.796184054
$ws.Range("Z9").Value = "17:58"
$ws.Range("AB9").Value = "17:58"

# Row 10
$ws.Range("A10").Value = 101775564
$ws.Range("I10").Value = "'12"
$ws.Range("J10").Value = "plantor/tuvor"
$ws.Range("Q10").Value = 567920.0505219861
$ws.Range("R10").Value = 6676823.397400577
$ws.Range("AW10").Value = "FREDRIK  Månsson "
$ws.Range("AX10").Value = "FREDRIK  Månsson "

# Row 11
$ws.Range("A11").Value = 101775565
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("Q11").Value = 567925.7637084418
$ws.Range("R11").Value = 6676837.41890596
$ws.Range("Z11").Value = "17:41"
$ws.Range("AB11").Value = "17:41"
$ws.Range("AW11").Value = "fanny westling"
$ws.Range("AX11").Value = "fanny westling"

# Row 12
$ws.Range("A12").Value = 101779916
$ws.Range("I12").Value = "'6"
$ws.Range("K12").Value = "blomknopp"
$ws.Range("Q12").Value = 567972.6507871181
$ws.Range("R12").Value = 6676859.163185491
$ws.Range("Z12").Value = "20:43"
$ws.Range("AB12").Value = "20:43"

# Row 13
$ws.Range("A13").Value = 101775890
$ws.Range("I13").Value = "'4"
$ws.Range("Q13").Value = 567873.0981311289
$ws.Range("R13").Value = 6676805.131272291
$ws.Range("Z13").Value = "18:01"
$ws.Range("AB13").Value = "18:01"

# Row 14
$ws.Range("A14").Value = 101775738
$ws.Range("I14").Value = "'50"
$ws.Range("J14").Value = "plantor/tuvor"
$ws.Range("Q14").Value = 567873.9455013226
$ws.Range("R14").Value = 6676839.93550191
$ws.Range("Z14").Value = "17:53"
$ws.Range("AB14").Value = "17:53"

# Row 15
$ws.Range("A15").Value = 101775684
$ws.Range("Q15").Value = 567898.9561973718
$ws.Range("R15").Value = 6676833.442172489
$ws.Range("Z15").Value = "17:49"
$ws.Range("AB15").Value = "17:49"

# Row 16
$ws.Range("A16").Value = 101774950
$ws.Range("K16").ClearContents()
$ws.Range("Q16").Value = 567968.1916468774
$ws.Range("R16").Value = 6676831.24963467
$ws.Range("Z16").Value = "17:25"
$ws.Range("AB16").Value = "17:25"

# Row 17
$ws.Range("A17").Value = 101775715
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()
$ws.Range("K17").Value = "blomknopp"
$ws.Range("Q17").Value = 567878.7189128834
$ws.Range("R17").Value = 6676850.957572321
$ws.Range("Z17").Value = "17:52"
$ws.Range("AB17").Value = "17:52"

# Row 18
$ws.Range("A18").Value = 101775512
$ws.Range("I18").Value = "'5"
$ws.Range("Q18").Value = 567927.6988580292
$ws.Range("R18").Value = 6676840.436715043
$ws.Range("Z18").Value = "17:38"
$ws.Range("AB18").Value = "17:38"

# Row 19
$ws.Range("A19").Value = 101775339
$ws.Range("I19").Value = "'40"
$ws.Range("Q19").Value = 567961.058589388
$ws.Range("R19").Value = 6676840.062643005
$ws.Range("Z19").Value = "17:31"
$ws.Range("AB19").Value = "17:31"

# Row 20
$ws.Range("A20").Value = 101775703
$ws.Range("B20").Value = 96334
$ws.Range("D20").Value = "VU"
$ws.Range("E20").Value = 220787
$ws.Range("F20").Value = "Knärot"
$ws.Range("G20").Value = "Goodyera repens"
$ws.Range("H20").Value = "(L.) R. Br."
$ws.Range("I20").Value = "'30"
$ws.Range("J20").Value = "plantor/tuvor"
$ws.Range("Q20").Value = 567893.7953871277
$ws.Range("R20").Value = 6676843.285874985
$ws.Range("Z20").Value = "17:51"
$ws.Range("AB20").Value = "17:51"

